$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45: date moved from 7/3/2022 to 7/4/2022 (entries continue past midnight)
$ws.Range("B45").Value = 44746

# New row 46: copy formatting from row 45, then set the actual values/formula
$ws.Range("A45:G45").Copy() | Out-Null
$ws.Range("A46:G46").PasteSpecial(-4122) | Out-Null

$ws.Range("A46").Value = 45
$ws.Range("B46").Value = 44746
$ws.Range("C46").Value = 0.33333333333333331
$ws.Range("D46").Value = 0.35416666666666669
$ws.Range("E46").Formula = "=D46-C46"
$ws.Range("F46").Value = "Code"
$ws.Range("G46").Value = "1. Unet model train for 12ep on 360 640 dataset`n2. FCN r50 model train for 12ep on 360 640 dataset"

# Row height for the wrapped two-line description (matches the other
# multi-line rows in the sheet, which are all 15pt * line-count)
$ws.Range("A46:G46").RowHeight = 30

# Selection moves to G45 as recorded by the author's last click before saving
$ws.Range("G45").Select() | Out-Null
